# Updates cryptos list prices/volumes (and reorders the Avalanche/Dogecoin
# and Hedera/WEMIXToken rows) to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.668.69'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.466.14'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.506'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0850'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').Value = '2.845.89'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.33%  '
$ws.Range('D16').Value = '2.468.97'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.787'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').Value = '41.581.72'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = '0.0₃0938'
$ws.Range('E20').Value = '  -2.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.50%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0764'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.90'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.22%  '
$ws.Range('D43').Value = '2.000.54'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.47'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.76%  '
$ws.Range('D48').Value = '2.725.72'
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '96.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '66.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.74%  '
